# Proceeding to: preprocess input,... > create entities
#
# 1) Refresh the cached "datetimeFigureOut" date field (Insert > Header &
#    Footer > Date and time > Update automatically > Apply to All) from
#    8/14/19 to 8/16/19 on the slide master and every slide layout.
# 2) Widen the labels on the two "Snip Diagonal Corner" example shapes on
#    slide 1: "Example 1" -> "Example 1: wide", "Example 2" -> "Example 2: wide".

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
            if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$newDate = "8/16/19"

# Slide master
$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes $newDate

# Every custom (slide) layout hanging off the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# Slide 1 text updates
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $text = $shape.TextFrame.TextRange.Text
        if ($text -eq "Example 1") {
            $shape.TextFrame.TextRange.Text = "Example 1: wide"
        } elseif ($text -eq "Example 2") {
            $shape.TextFrame.TextRange.Text = "Example 2: wide"
        }
    }
}
